# Update the ARCHITECTURE sheet's "type_wall" (column O) values to use the
# new SG wall types "T8" (was "T2") and "T7" (was "T6"), based on the BCA
# ETTV example, and leave the sheet selection positioned on the last edited
# cell (O20), matching the author's final selection/active-sheet state.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ARCHITECTURE")

# Rows that previously referenced wall type "T2" now use "T8"
$t8Rows = @(2, 3, 14, 16)
foreach ($r in $t8Rows) {
    $ws.Range("O$r").Value = "T8"
}

# Rows that previously referenced wall type "T6" now use "T7"
$t7Rows = @(4, 5, 6, 7, 8, 9, 10, 11, 12, 13, 15, 17, 18, 19, 20)
foreach ($r in $t7Rows) {
    $ws.Range("O$r").Value = "T7"
}

# Leave the ARCHITECTURE sheet active/selected with the bottom-right pane's
# selection on O20 (the last cell touched), and INTERNAL_LOADS no longer
# tabSelected.
$ws.Range("O20").Select() | Out-Null
